$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Ingreso" (sheet1): append new aporte rows 712-720
# ---------------------------------------------------------------------------
$wsIngreso = $wb.Worksheets.Item("Ingreso")
$wsIngreso.Activate() | Out-Null

$ingresoRows = @(
    @(712, 45417, "Michy",     100),
    @(713, 45417, "Anuel",     100),
    @(714, 45438, "Invitados", 300),
    @(715, 45438, "Julio",     100),
    @(716, 45438, "Yeyo",       80),
    @(717, 45438, "Rubio",     300),
    @(718, 45438, "Punto",     200),
    @(719, 45438, "Johan",     300),
    @(720, 45438, "Anuel",     100)
)

foreach ($row in $ingresoRows) {
    $r = $row[0]
    $wsIngreso.Cells.Item($r, 1).Value = $row[1]
    $wsIngreso.Cells.Item($r, 2).Value = $row[2]
    $wsIngreso.Cells.Item($r, 3).Value = $row[3]
    $wsIngreso.Cells.Item($r, 4).Value = "Aporte"
}

$wsIngreso.Range("A720").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Gastos" (sheet2): append new gasto rows 91-92
# ---------------------------------------------------------------------------
$wsGastos = $wb.Worksheets.Item("Gastos")
$wsGastos.Activate() | Out-Null

$gastosRows = @(
    @(91, 45417, "Arbitro y agua", 960),
    @(92, 45438, "Agua",           160)
)

foreach ($row in $gastosRows) {
    $r = $row[0]
    $wsGastos.Cells.Item($r, 1).Value = $row[1]
    $wsGastos.Cells.Item($r, 2).Value = $row[2]
    $wsGastos.Cells.Item($r, 3).Value = $row[3]
}

$wsGastos.Range("A92").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Histórico de tecnicas" (sheet4): view scrolled down slightly; no
# data changes. Keep the existing selection on A14.
# ---------------------------------------------------------------------------
$wsHistorico = $wb.Worksheets.Item("Histórico de tecnicas")
$wsHistorico.Activate() | Out-Null
$wsHistorico.Range("A14").Select() | Out-Null

# ---------------------------------------------------------------------------
# Restore "Ingreso" as the active sheet/tab, matching the original workbook.
# ---------------------------------------------------------------------------
$wsIngreso.Activate() | Out-Null
$wsIngreso.Range("A720").Select() | Out-Null
